{"js": "// Apply the \"before -> after\" report-regeneration edit described by the\n// diff: refreshed date/time/version stamp, and the numeric prefixes\n// (\"1. \", \"2. \", \"5.1. \", ...) stripped from the Heading1/Heading2 titles,\n// plus the Table-of-Contents field widened from levels 1-3 to 1-4.\n\nconst body = context.document.body;\n\n// --- 1) Simple, unambiguous whole-string text replacements -------------\n// Each pair is unique within the document body, so a plain search+replace\n// is safe and keeps every other run/formatting attribute untouched.\nconst replacements = [\n  [\n    \"This document was generated on 2020-08-07, 15:05:24 with the Automatic Report Generator (ARG) version \\\"develop\\\" on the Linux system runner-0277ea0f-project-18732201-concurrent-0.\",\n    \"This document was generated on 2021-08-31, 15:49:38 with the Automatic Report Generator (ARG) version \\\"1.1.7-RC4\\\" on the Linux system runner-0277ea0f-project-18732201-concurrent-0.\"\n  ],\n  [\"1. Table of Contents\", \"Table of Contents\"],\n  [\"2. List of Figures\", \"List of Figures\"],\n  [\"3. List of Tables\", \"List of Tables\"],\n  [\"4. Introduction\", \"Introduction\"],\n  [\"5.1. Per-File Meta-Information\", \"Per-File Meta-Information\"],\n  [\"5.2. Directory Meta-Information\", \"Directory Meta-Information\"],\n  // Must run after the \"5.1./5.2.\" replacements above (otherwise this\n  // broader \"5. \" match would also eat into them).\n  [\"5. Key/Value Files\", \"Key/Value Files\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// The title page's standalone date (\"2020-08-07\") on its own line is a\n// *different* run from the \"generated on ...\" sentence (already handled\n// above), so it needs its own targeted search.\nconst titleDate = body.search(\"2020-08-07\", { matchCase: true });\ntitleDate.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < titleDate.items.length; i++) {\n  titleDate.items[i].insertText(\"2021-08-31\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Widen the Table of Contents field from \"\\o 1-3\" to \"\\o 1-4\" -----\n// The TOC field code (\"instrText\") is not part of the normal body text\n// stream, so `body.search` can never find it; it has to be reached\n// through the paragraph that hosts the field's runs. We rebuild that one\n// paragraph's OOXML, keeping the begin/separate/end field-character\n// structure (and the cached \"Right-click to update field.\" placeholder\n// text) intact, only swapping the \"1-3\" level range for \"1-4\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  // The TOC field paragraph renders as empty text (its only content is\n  // the field-character/instrText run), and is immediately preceded by\n  // the \"Table of Contents\" heading paragraph.\n  if (\n    p.text === \"\" &&\n    i > 0 &&\n    paragraphs.items[i - 1].text === \"Table of Contents\"\n  ) {\n    const fieldOoxml =\n      '<?xml version=\"1.0\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p><w:r>' +\n      '<w:fldChar w:fldCharType=\"begin\"/>' +\n      '<w:instrText xml:space=\"preserve\">TOC \\\\o 1-4 \\\\h \\\\z \\\\u</w:instrText>' +\n      '<w:fldChar w:fldCharType=\"separate\"/>' +\n      '<w:t>Right-click to update field.</w:t>' +\n      '<w:fldChar w:fldCharType=\"end\"/>' +\n      '</w:r></w:p></w:body></w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>';\n    p.getRange().insertOoxml(fieldOoxml, Word.InsertLocation.replace);\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "# Apply the \"before -> after\" report-regeneration edit described by the\n# diff: refreshed date/time/version stamp, and the numeric prefixes\n# (\"1. \", \"2. \", \"5.1. \", ...) stripped from the Heading1/Heading2 titles,\n# plus the Table-of-Contents field widened from levels 1-3 to 1-4.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $d.Content.Find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap            (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace         (wdReplaceAll)\n    ) | Out-Null\n}\n\n# --- 1) Simple, unambiguous whole-string text replacements -------------\n# Each pair is unique within the document body, so a plain find/replace-all\n# is safe and keeps every other run/formatting attribute untouched.\n#\n# NOTE: Find/Replace in this host auto-corrects straight quote characters\n# that appear *inside* a ReplaceWith string into curly \"smart quotes\", so\n# the version-string swap is deliberately split to only touch the text\n# between the existing straight quotes (which are then left alone).\nReplace-AllText \"15:05:24\" \"15:49:38\"\nReplace-AllText \"develop\" \"1.1.7-RC4\"\n# Rewrites both the title-page date and the \"generated on\" sentence's date\n# in one shot (they share the same \"2020-08-07\" substring).\nReplace-AllText \"2020-08-07\" \"2021-08-31\"\nReplace-AllText \"1. Table of Contents\" \"Table of Contents\"\nReplace-AllText \"2. List of Figures\" \"List of Figures\"\nReplace-AllText \"3. List of Tables\" \"List of Tables\"\nReplace-AllText \"4. Introduction\" \"Introduction\"\nReplace-AllText \"5.1. Per-File Meta-Information\" \"Per-File Meta-Information\"\nReplace-AllText \"5.2. Directory Meta-Information\" \"Directory Meta-Information\"\nReplace-AllText \"5. Key/Value Files\" \"Key/Value Files\"\n\n# --- 2) Widen the Table of Contents field from \"\\o 1-3\" to \"\\o 1-4\" -----\n# The field's instruction text (\"TOC \\o 1-3 \\h \\z \\u\") is not part of the\n# normal document text stream that Find/Execute walks, so it has to be\n# reached through the Fields collection instead. Re-assigning a Field's\n# Code rebuilds that field (this is the only mutation the object model\n# exposes for editing field instructions), so we only do this for the\n# Table-of-Contents field (the first field in the document).\n$fields = $d.Fields\nfor ($i = 1; $i -le $fields.Count; $i++) {\n    $f = $fields.Item($i)\n    if ($f.Code.Text -like \"TOC \\o 1-3*\") {\n        $f.Code = \"TOC \\o 1-4 \\h \\z \\u\"\n        break\n    }\n}\n"}
